$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B, C, D, F, G, H, J, K, L, M, N
$colIndices = @(2, 3, 4, 6, 7, 8, 10, 11, 12, 13, 14)

# New loading_percent values per row (row 2..25), one array per row
# ordered to match $colIndices above
$data = @(
    @(16.9518374939988, 7.310502714159211, 7.15717425555001, 43.95629878570092, 53.02591298662058, 20.42477573345007, 10.91876868098595, 12.93339559677725, 11.389587354356, 17.11304594603876, 23.49642253661797),
    @(16.82118491008674, 7.254373317559877, 7.154887667678968, 44.01406064084343, 53.07293093642004, 20.46839479683171, 10.93646102320808, 12.83854812474575, 11.40136862263742, 17.10670021546743, 23.55762513097128),
    @(16.74416312562985, 7.219025625961913, 7.154281768719701, 44.05742211633562, 53.11296171954631, 20.49800635360531, 10.94795671699744, 12.78267821209051, 11.40986157194969, 17.10532315283803, 23.59708526932441),
    @(16.71361102525151, 7.204402828759903, 7.15423650389187, 44.07707630547512, 53.13207824122482, 20.51078467470372, 10.95280078888767, 12.76052729585183, 11.41363972109811, 17.10539780143324, 23.61363983958341),
    @(16.70858915699459, 7.201961654745421, 7.154241193473602, 44.0804596513484, 53.13542175044478, 20.51294946751162, 10.95361478840068, 12.75688699358462, 11.4142862541414, 17.105448661413, 23.61641738970416),
    @(16.74374766980256, 7.218829296524421, 7.154280340577063, 44.05767914846599, 53.11320818456746, 20.49817580633611, 10.94802139953072, 12.78237695335006, 11.40991124031762, 17.10532158230668, 23.59730660837285),
    @(16.906143106292, 7.291335136469404, 7.156220791717804, 43.97457568877353, 53.03980667221757, 20.43922833552049, 10.92473797914875, 12.90021408940772, 11.39338853617247, 17.11033653741783, 23.51713536604868),
    @(17.24858521739947, 7.426335183165052, 7.166314554260489, 43.87430271960156, 52.98452515719476, 20.34608445585417, 10.88407854756522, 13.14908603196569, 11.3709533580675, 17.14004892379439, 23.37479650671814),
    @(17.51284170944031, 7.52092069518977, 7.177502950456045, 43.83889231431975, 52.99803375420377, 20.29134228090491, 10.85722661134363, 13.34140057073004, 11.36051068019294, 17.17383520995588, 23.2792153099677),
    @(17.635373431472, 7.562906131953513, 7.183397798995291, 43.83108982615311, 53.01592478128612, 20.26941153078736, 10.84566110339523, 13.4306379980296, 11.35706356251688, 17.19176202618065, 23.23767046915844),
    @(17.68206920367139, 7.578651265885894, 7.185744489017576, 43.82932846515087, 53.02438590939599, 20.2615342039679, 10.84137451880451, 13.4646553986648, 11.35594488713342, 17.19891433462835, 23.22221570849771),
    @(17.67199989119472, 7.575267174530261, 7.185234019733973, 43.82965475820595, 53.02248871918052, 20.26321171959528, 10.84229358040083, 13.45731956125573, 11.35617752377231, 17.19735784160814, 23.22553184674482),
    @(17.63920941360308, 7.564204599421194, 7.183588577002411, 43.83092101129373, 53.0165871104966, 20.26875489312776, 10.84530658139078, 13.43343227319953, 11.35696779159476, 17.19234318295266, 23.2363934432406),
    @(17.61916167220448, 7.557408300816418, 7.182595558117328, 43.83185198477378, 53.01319168988273, 20.27220590385704, 10.84716423436496, 13.41882913460009, 11.35747614115568, 17.18931881797853, 23.24308257818799),
    @(17.50487742245175, 7.518155536698663, 7.177133797030693, 43.83956927926884, 52.99710067440687, 20.2928353104356, 10.85799548250957, 13.33560166873379, 11.36076212091421, 17.17271478008978, 23.28196922855071),
    @(17.4353350783133, 7.493805486265546, 7.173988473647759, 43.84643031715377, 52.99023611707434, 20.30625195969284, 10.86480619805333, 13.28497415079621, 11.36311125712596, 17.16318101342485, 23.30631997745147),
    @(17.39555684162857, 7.479702139216267, 7.172255239734302, 43.85115849481389, 52.98739363337989, 20.31424857608777, 10.86878470514137, 13.25602113471588, 11.36458513856743, 17.15793839635779, 23.32050812112985),
    @(17.3821276672699, 7.474910307018535, 7.171681467772531, 43.8528936894466, 52.98662121146873, 20.31700413232311, 10.87014227600457, 13.24624756821642, 11.36510526877595, 17.15620483534518, 23.32534331160938),
    @(17.44271541392922, 7.496407746248559, 7.174315455698585, 43.84561903145015, 52.99085242864518, 20.30479478461278, 10.86407485835309, 13.29034648867574, 11.36284849170577, 17.16417098893982, 23.30370894478306),
    @(17.64883305135977, 7.567458152565258, 7.184068788619023, 43.83051670858436, 53.0182748223501, 20.26711512984054, 10.84441906841956, 13.44044265455276, 11.35673061062269, 17.19380626854511, 23.23319560889911),
    @(17.78524845171757, 7.612995205348616, 7.191109524720566, 43.82760086013626, 53.04602311949371, 20.24498033584886, 10.83211490218921, 13.53983874266972, 11.35381992439889, 17.21529320900535, 23.18872735272388),
    @(17.71229752922694, 7.588774737148844, 7.18729123920384, 43.82852128588976, 53.03031549816891, 20.25656616057739, 10.83863239798452, 13.48667921385743, 11.35527414263665, 17.20363273797862, 23.21231332004133),
    @(17.43937813027391, 7.49523158872704, 7.174167393240694, 43.84598337249675, 52.99057035490134, 20.30545269127058, 10.86440530096136, 13.28791717478354, 11.36296690373362, 17.16372267789263, 23.30488880515167),
    @(17.15357595195471, 7.390606043746146, 7.162916760861552, 43.89471127669635, 52.98997646602677, 20.36887821391398, 10.89454561379961, 13.07999412765566, 11.37595922740056, 17.12989847731886, 23.41171799879289)
)

$startRow = 2
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $colIndices.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $colIndices[$c]).Value = $rowValues[$c]
    }
}
